# Generate Report for Handback
# Updates the localization-status workbook to reflect that both zh-cn and
# de-de targets have been handed back (now in sync with en-US), and fills
# in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the per-locale detail sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status cells and widen their
# columns so the longer status text is readable.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Helper data: per-locale target / handback file names and handback time
# ---------------------------------------------------------------------
$rows = @(
    @{
        SourceDisplay = "780cf0fa-05e7-446e-ba54-f07279d0176f.md"
        SourceAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef3321d2e0a2c20ce7c2ca589c7130bc18792f90/e2e/780cf0fa-05e7-446e-ba54-f07279d0176f.md"
        Row = 2
        ZhHandback = "780cf0fa-05e7-446e-ba54-f07279d0176f.abc6b9ee71645279fa20319e3fd92ea04ec1dd93.zh-cn.xlf"
        DeHandback = "780cf0fa-05e7-446e-ba54-f07279d0176f.abc6b9ee71645279fa20319e3fd92ea04ec1dd93.de-de.xlf"
    },
    @{
        SourceDisplay = "d7be68c9-10c8-40c6-8f3a-7e6587dd5b6e.md"
        SourceAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef3321d2e0a2c20ce7c2ca589c7130bc18792f90/e2e/d7be68c9-10c8-40c6-8f3a-7e6587dd5b6e.md"
        Row = 3
        ZhHandback = "d7be68c9-10c8-40c6-8f3a-7e6587dd5b6e.993665eb2dae399ba795baef9131f96704058da7.zh-cn.xlf"
        DeHandback = "d7be68c9-10c8-40c6-8f3a-7e6587dd5b6e.993665eb2dae399ba795baef9131f96704058da7.de-de.xlf"
    }
)

$zhHandbackDateTime = "2016-09-01 18:28:17"
$deHandbackDateTime = "2016-09-01 18:28:25"

# ---------------------------------------------------------------------
# zh-cn sheet: "Latest Target File" (I), "Latest Handback File" (J) and
# "Latest Handback DateTime" (K) columns.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

foreach ($r in $rows) {
    $targetCell = $zh.Cells.Item($r.Row, 9)
    $targetCell.Value = $r.SourceDisplay
    $zh.Hyperlinks.Add($targetCell, $r.SourceAddress, "", "", $r.SourceDisplay) | Out-Null

    $zh.Cells.Item($r.Row, 10).Value = $r.ZhHandback
    $zh.Cells.Item($r.Row, 11).Value = $zhHandbackDateTime
}

# ---------------------------------------------------------------------
# de-de sheet: same three columns, with the de-de handoff/handback data.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

foreach ($r in $rows) {
    $targetCell = $de.Cells.Item($r.Row, 9)
    $targetCell.Value = $r.SourceDisplay
    $de.Hyperlinks.Add($targetCell, $r.SourceAddress, "", "", $r.SourceDisplay) | Out-Null

    $de.Cells.Item($r.Row, 10).Value = $r.DeHandback
    $de.Cells.Item($r.Row, 11).Value = $deHandbackDateTime
}
